$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "E1"  = 1.399822193902046
    "F1"  = -1.570796365408233
    "E2"  = 1.398641107009634
    "F2"  = -1.570796365220714
    "E3"  = 1.391275308874927
    "F3"  = -1.570796364051259
    "E4"  = 1.373855546651759
    "F4"  = -1.570796361285554
    "E5"  = 1.34470615965872
    "F5"  = -1.570796356657557
    "E6"  = 1.303979480685029
    "F6"  = -1.570796350191454
    "E7"  = 1.25329023729641
    "F7"  = -1.570796342143612
    "E8"  = 1.195349953140961
    "F8"  = -1.570796332944535
    "E9"  = 1.133601349255035
    "F9"  = -1.570796323140818
    "E10" = 1.071852745369108
    "F10" = -1.570796313337101
    "E11" = 1.013912461213659
    "F11" = -1.570796304138024
    "E12" = 0.9632232178250399
    "F12" = -1.570796296090182
    "E13" = 0.9224965388513493
    "F13" = -1.570796289624078
    "E14" = 0.8933471518583104
    "F14" = -1.570796284996081
    "E15" = 0.8759273896351421
    "F15" = -1.570796282230376
    "E16" = 0.8685615915004347
    "F16" = -1.570796281060922
    "E17" = 0.8673805046080224
    "F17" = -1.570796280873402
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
